$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A "header" labels entered first, in the order the author typed them
# (this determines the shared-string table ordering seen in the target file)
$ws.Range("A19").Value = "Unclosed:"
$ws.Range("B19").Value = 'Unclosed "double quote'
$ws.Range("A20").Value = "BigIntegers:"
$ws.Range("A24").Value = "BigDecimals:"

# BigInteger / BigDecimal expression cells
$ws.Range("B20").Value = '${biZero}'
$ws.Range("B21").Value = '${biAnswer}'
$ws.Range("B22").Value = '${biBiggerThanLong}'
$ws.Range("B23").Value = '${biBiggerThanDouble}'
$ws.Range("B24").Value = '${bdZero}'
$ws.Range("B25").Value = '${bdAnswer}'
$ws.Range("B26").Value = '${bdSmallerThanNormal}'
$ws.Range("B27").Value = '${bdBiggerThanDouble}'

# ValueHolder block
$ws.Range("A28").Value = "ValueHolder:"
$ws.Range("B28").Value = '${valueHolder.answer}'
$ws.Range("B29").Value = '${valueHolder.IHaveAQuestion}'

# New column width for column A, matching the <cols> entry added in the diff
# (target stored width is 12.42578125 characters; closest reachable value via
# the ColumnWidth COM property's pixel-grid rounding is used here)
$ws.Columns.Item(1).ColumnWidth = 11.6
